# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# The two sheets list mostly the same events; "全部类型" simply has one
# extra row near the top (a 演出 entry), so its data rows are shifted
# down by one relative to "展览".

$wb = $excel.ActiveWorkbook

$updates_exhibition = @{
    "F3"  = 658    # 合肥·运动番only·群青日和            659 -> 658
    "F4"  = 243    # 合肥·FT动漫嘉年华（免费）             242 -> 243
    "F6"  = 10006  # 合肥·第六届环形宇宙动漫游戏嘉年华... 9996 -> 10006
    "F7"  = 905    # 合肥·第二届华盟动漫次元嘉年华         904 -> 905
    "F9"  = 1245   # 合肥·Look Look动漫嘉年华             1244 -> 1245
    "F10" = 5527   # 合肥·城市动漫节                     5525 -> 5527
    "F15" = 71     # 合肥·星光次元动漫文化节              70 -> 71
}

$updates_all_types = @{
    "F4"  = 658    # 合肥·运动番only·群青日和            659 -> 658
    "F5"  = 243    # 合肥·FT动漫嘉年华（免费）             242 -> 243
    "F7"  = 10006  # 合肥·第六届环形宇宙动漫游戏嘉年华... 9996 -> 10006
    "F8"  = 905    # 合肥·第二届华盟动漫次元嘉年华         904 -> 905
    "F10" = 1245   # 合肥·Look Look动漫嘉年华             1244 -> 1245
    "F11" = 5527   # 合肥·城市动漫节                     5525 -> 5527
    "F16" = 71     # 合肥·星光次元动漫文化节              70 -> 71
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($cell in $updates_exhibition.Keys) {
    $wsExhibition.Range($cell).Value = $updates_exhibition[$cell]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($cell in $updates_all_types.Keys) {
    $wsAllTypes.Range($cell).Value = $updates_all_types[$cell]
}
